$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1566.9286
$ws.Range("I4").Value = 497.6
$ws.Range("J4").Value = 2161
$ws.Range("K4").Value = 497.6
$ws.Range("L4").Value = 2161
$ws.Range("M4").Value = -383.6
$ws.Range("N4").Value = -2389
$ws.Range("H103").Value = 2001.25
$ws.Range("I103").Value = 1000
$ws.Range("J103").Value = 2335
$ws.Range("K103").Value = 3000
$ws.Range("L103").Value = 7005
$ws.Range("M103").Value = -2414
$ws.Range("N103").Value = -8177
$ws.Range("H107").Value = 601.7857
$ws.Range("I107").Value = 601.7857
$ws.Range("K107").Value = 601.7857
$ws.Range("M107").Value = 1318.2143
$ws.Range("H112").Value = 5876.625
$ws.Range("J112").Value = 1595.0667
$ws.Range("L112").Value = 4785.2001
$ws.Range("N112").Value = -7001.2001
$ws.Range("H129").Value = 973.66
$ws.Range("I129").Value = 429.16666
$ws.Range("J129").Value = 1047.909
$ws.Range("K129").Value = 1287.49998
$ws.Range("L129").Value = 3143.727
$ws.Range("M129").Value = 3712.50002
$ws.Range("N129").Value = -13143.727
$ws.Range("H137").Value = 1044638.75
$ws.Range("I137").Value = 3447.5881
$ws.Range("J137").Value = 2224655.5
$ws.Range("K137").Value = 10342.7643
$ws.Range("L137").Value = 6673966.5
$ws.Range("M137").Value = -7792.764299999999
$ws.Range("N137").Value = -6679066.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H96").Value = 27000
$ws.Range("J96").Value = 27000
$ws.Range("L96").Value = 27000
$ws.Range("N96").Value = -32492
$ws.Range("H122").Value = 13889832
$ws.Range("I122").Value = 1061.25
$ws.Range("J122").Value = 125000000
$ws.Range("K122").Value = 3183.75
$ws.Range("L122").Value = 375000000
$ws.Range("M122").Value = -733.75
$ws.Range("N122").Value = -375004900
$ws.Range("H129").Value = 30246.334
$ws.Range("J129").Value = 30246.334
$ws.Range("L129").Value = 30246.334
$ws.Range("N129").Value = -40246.334
$ws.Range("H137").Value = 56597.5
$ws.Range("J137").Value = 56597.5
$ws.Range("L137").Value = 56597.5
$ws.Range("N137").Value = -66797.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2456.25
$ws.Range("I107").Value = 2169.2307
$ws.Range("J107").Value = 3700
$ws.Range("K107").Value = 2169.2307
$ws.Range("L107").Value = 3700
$ws.Range("M107").Value = -249.2307000000001
$ws.Range("N107").Value = -7540
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2598.8
$ws.Range("I86").Value = 2700
$ws.Range("K86").Value = 2700
$ws.Range("M86").Value = -1577
$ws.Range("H89").Value = 2598.8
$ws.Range("I89").Value = 2700
$ws.Range("K89").Value = 13500
$ws.Range("M89").Value = -7884
$ws.Range("H94").Value = 807.2
$ws.Range("J94").Value = 837.3333
$ws.Range("L94").Value = 837.3333
$ws.Range("N94").Value = -1739.3333
$ws.Range("H123").Value = 59400
$ws.Range("J123").Value = 59400
$ws.Range("L123").Value = 59400
$ws.Range("N123").Value = -69200

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2461.7917
$ws.Range("I68").Value = 787.45
$ws.Range("J68").Value = 4554.7188
$ws.Range("K68").Value = 2362.35
$ws.Range("L68").Value = 13664.1564
$ws.Range("M68").Value = -1551.35
$ws.Range("N68").Value = -15286.1564
$ws.Range("H71").Value = 2461.7917
$ws.Range("I71").Value = 787.45
$ws.Range("J71").Value = 4554.7188
$ws.Range("K71").Value = 7087.05
$ws.Range("L71").Value = 40992.4692
$ws.Range("M71").Value = -3031.05
$ws.Range("N71").Value = -49104.4692
$ws.Range("H93").Value = 4921.1665
$ws.Range("J93").Value = 4921.1665
$ws.Range("L93").Value = 14763.4995
$ws.Range("N93").Value = -18507.4995
$ws.Range("H131").Value = 24997.6
$ws.Range("I131").Value = 1267.1428
$ws.Range("J131").Value = 30031.334
$ws.Range("K131").Value = 3801.4284
$ws.Range("L131").Value = 90094.00199999999
$ws.Range("M131").Value = 1238.5716
$ws.Range("N131").Value = -100174.002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2272.7273
$ws.Range("I113").Value = 2175
$ws.Range("J113").Value = 2533.3333
$ws.Range("K113").Value = 2175
$ws.Range("L113").Value = 2533.3333
$ws.Range("M113").Value = -5
$ws.Range("N113").Value = -6873.3333
$ws.Range("H122").Value = 7733.25
$ws.Range("I122").Value = 8879.799999999999
$ws.Range("J122").Value = 2000.5
$ws.Range("K122").Value = 26639.4
$ws.Range("L122").Value = 6001.5
$ws.Range("M122").Value = -24189.4
$ws.Range("N122").Value = -10901.5
$ws.Range("H126").Value = 3103.8462
$ws.Range("I126").Value = 1994.4445
$ws.Range("J126").Value = 5600
$ws.Range("K126").Value = 5983.333500000001
$ws.Range("L126").Value = 16800
$ws.Range("M126").Value = -3513.333500000001
$ws.Range("N126").Value = -21740
$ws.Range("H134").Value = 34215.777
$ws.Range("J134").Value = 34215.777
$ws.Range("L134").Value = 102647.331
$ws.Range("N134").Value = -107717.331
$ws.Range("H137").Value = 49800
$ws.Range("J137").Value = 49800
$ws.Range("L137").Value = 49800
$ws.Range("N137").Value = -60000

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 12800
$ws.Range("I5").Value = 12800
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 12800
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -12687
$ws.Range("N5").ClearContents()
$ws.Range("H7").Value = 3000
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224
$ws.Range("H22").Value = 10000
$ws.Range("J22").Value = 10000
$ws.Range("L22").Value = 10000
$ws.Range("N22").Value = -10590
$ws.Range("H27").Value = 10000
$ws.Range("J27").Value = 10000
$ws.Range("L27").Value = 10000
$ws.Range("N27").Value = -10214
$ws.Range("H40").Value = 4427.857
$ws.Range("I40").Value = 3498.3333
$ws.Range("J40").Value = 10005
$ws.Range("K40").Value = 3498.3333
$ws.Range("L40").Value = 10005
$ws.Range("M40").Value = -3362.3333
$ws.Range("N40").Value = -10277
$ws.Range("H122").Value = 7509.091
$ws.Range("I122").Value = 7260
$ws.Range("J122").Value = 10000
$ws.Range("K122").Value = 21780
$ws.Range("L122").Value = 30000
$ws.Range("M122").Value = -19330
$ws.Range("N122").Value = -34900
$ws.Range("H126").Value = 3000
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H109").Value = 78800
$ws.Range("J109").Value = 78800
$ws.Range("L109").Value = 78800
$ws.Range("N109").Value = -81574
